{"js": "// Office.js (Word JavaScript API) script.\n// Applies the \"v2\" update to the instructor_materials README table:\n//  1. Widen the table slightly and re-distribute the 4 column widths.\n//  2. Rename the \"Shiny_App_Student_Questions_FINAL.pdf\" row to\n//     \"Student_Handout_FINAL.pdf\" and give it the description text that used\n//     to live on the (now removed) \"Pre-module_ Student_Handout_FINAL.pdf\"\n//     row.\n//  3. Delete the old \"Pre-module_ Student_Handout_FINAL.pdf\" row entirely\n//     (its content was folded into the row above).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// --- 1. Table width + column widths (values are in points; OOXML stores\n//     twips/dxa = points * 20) -------------------------------------------\n// table.width has no real setter in this host, so use the lower-level OM\n// bridge (Table.PreferredWidth, same thing COM/VBA would set) instead.\ntable._omSet(\"PreferredWidth\", 9430 / 20, \"Table\");\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nrows.items[0].cells.load(\"items\");\nawait context.sync();\n\nconst newColumnWidthsDxa = [2410, 1440, 1800, 3780];\nconst headerCells = rows.items[0].cells.items;\nfor (let col = 0; col < newColumnWidthsDxa.length; col++) {\n  // TableCell.columnWidth resizes the whole column (gridCol + every cell's\n  // tcW), matching what the diff shows happening on every row.\n  headerCells[col].columnWidth = newColumnWidthsDxa[col] / 20;\n}\nawait context.sync();\n\n// --- 2 & 3. Row content -----------------------------------------------\nrows.items[4].cells.load(\"items\");\nrows.items[5].cells.load(\"items\");\nawait context.sync();\n\nconst row4Cells = rows.items[4].cells.items;\nrow4Cells[0].value = \"Student_Handout_FINAL.pdf\";\nrow4Cells[3].value =\n  \"Handout for students to work through prior to completing the module. \" +\n  \"While this version has been archived as a pdf file, we refer interested \" +\n  \"readers to http://www.module5.macrosystemseddie.org/ for editable \" +\n  \"Microsoft Word files. We note that some changes may be made to the \" +\n  \"files on the website as they are updated over time.\";\nawait context.sync();\n\n// Remove the old \"Pre-module_ Student_Handout_FINAL.pdf\" row (index 5) \u2014\n// its text was merged into row 4 above.\nrows.items[5].delete();\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the \"v2\" update to the instructor_materials README table:\n#  1. Widen the table slightly and re-distribute the 4 column widths.\n#  2. Rename the \"Shiny_App_Student_Questions_FINAL.pdf\" row to\n#     \"Student_Handout_FINAL.pdf\" and give it the description text that used\n#     to live on the (now removed) \"Pre-module_ Student_Handout_FINAL.pdf\"\n#     row.\n#  3. Delete the old \"Pre-module_ Student_Handout_FINAL.pdf\" row entirely\n#     (its content was folded into the row above).\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# --- 1. Table width + column widths (Word COM widths are in points; OOXML\n#     stores twips/dxa = points * 20) --------------------------------------\n$tbl.PreferredWidth = 9430 / 20\n\n$tbl.Columns.Item(1).Width = 2410 / 20\n$tbl.Columns.Item(2).Width = 1440 / 20\n$tbl.Columns.Item(3).Width = 1800 / 20\n$tbl.Columns.Item(4).Width = 3780 / 20\n\n# --- 2. Row content (row 5 = \"Shiny_App_Student_Questions_FINAL.pdf\") -----\n$tbl.Cell(5, 1).Range.Text = \"Student_Handout_FINAL.pdf\"\n$tbl.Cell(5, 4).Range.Text = \"Handout for students to work through prior to completing the module. While this version has been archived as a pdf file, we refer interested readers to http://www.module5.macrosystemseddie.org/ for editable Microsoft Word files. We note that some changes may be made to the files on the website as they are updated over time.\"\n\n# --- 3. Remove the old \"Pre-module_ Student_Handout_FINAL.pdf\" row (row 6)\n#     \u2014 its text was merged into row 5 above.\n$tbl.Rows.Item(6).Delete()\n"}
